$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts existing "Numbers"/"Letters" columns to B/C)
$ws.Columns.Item(1).Insert()

# Populate new "URL" column: header + category/subcategory/id pattern rows
$ws.Range("A1").Value = "URL"
$ws.Range("A2").Value = "/category"
$ws.Range("A3").Value = "/subcategory"
$ws.Range("A4").Value = "/id{x}-item{y}"

# Apply the new (red) font color to the whole table, matching the
# font added to the workbook's style table
$ws.Range("A1:C4").Font.Color = 1974729
$ws.Range("B5:C18").Font.Color = 1974729

# Restore the cursor/selection position recorded in the saved file
$ws.Range("J16").Select() | Out-Null
